$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") advances by one day (46072 -> 46073) on every data row.
$ws.Range("C2").Value = 46073
$ws.Range("C3").Value = 46073
$ws.Range("C4").Value = 46073
$ws.Range("C5").Value = 46073
$ws.Range("C6").Value = 46073
$ws.Range("C7").Value = 46073
$ws.Range("C8").Value = 46073
$ws.Range("C9").Value = 46073
$ws.Range("C10").Value = 46073
$ws.Range("C11").Value = 46073
$ws.Range("C12").Value = 46073

# Rows 4-11 get reshuffled (Beteckning / Datum / Area columns).
$ws.Range("A4").Value = "A 50762-2025"
$ws.Range("B4").Value = 45946
$ws.Range("G4").Value = 2.7

$ws.Range("A5").Value = "A 14516-2023"
$ws.Range("B5").Value = 45012.86600694444
$ws.Range("G5").Value = 0.4

$ws.Range("A6").Value = "A 4156-2023"
$ws.Range("B6").Value = 44953
$ws.Range("G6").Value = 1.5

$ws.Range("A7").Value = "A 8679-2026"
$ws.Range("B7").Value = 46066
$ws.Range("G7").Value = 2.1

$ws.Range("A8").Value = "A 26262-2024"
$ws.Range("B8").Value = 45468.66077546297
$ws.Range("G8").Value = 0.6

$ws.Range("A10").Value = "A 14517-2023"
$ws.Range("B10").Value = 45012
$ws.Range("G10").Value = 0.6

$ws.Range("A11").Value = "A 23798-2024"
$ws.Range("B11").Value = 45455.43208333333
$ws.Range("G11").Value = 1.3
